# Append two new days of COVID-19 data (2020-05-25 and 2020-05-26) to the
# "Tabela1" table on the single worksheet, extending the table/autofilter
# range and the sheet selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item(1)

# Grow the table by two rows (this also extends ref/autoFilter on the table).
$row76 = $lo.ListRows.Add()
$row77 = $lo.ListRows.Add()

# New data, one array per row: Date, Tested(all), Tested(daily), Positive(all),
# Positive(daily), Hospitalized, Intensive care, Discharged, Deaths(all), Deaths(daily)
$newData = @(
    @(43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
    @(43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0)
)

for ($i = 0; $i -lt 2; $i++) {
    $r = 76 + $i

    # Match the formatting already used by the preceding data rows: date
    # column uses the custom date format, the "Tested (all)" column uses
    # thousands separators, the rest are plain numbers (the default, so no
    # explicit NumberFormat needed there). Font mirrors the rest of the
    # table ("Calibri Light" 10pt, right aligned).
    $ws.Cells.Item($r, 1).NumberFormat = "d/\ m/\ yyyy;@"
    $ws.Cells.Item($r, 2).NumberFormat = "#,##0"

    $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 10))
    $rowRange.Font.Name = "Calibri Light"
    $rowRange.Font.Size = 10
    $rowRange.HorizontalAlignment = -4152  # xlRight

    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $newData[$i][$c - 1]
    }
}

# Match the author's final selection after entering the last row of data.
[void]$ws.Range("A77:J77").Select()
